# Science Catalog page tweaks:
#  - slide with SlideID 260 (position 1): nudge the "Harvests" legend
#    textbox + its icon picture a bit to the left, and add a new
#    "Products" caption textbox near the other legend captions.
#  - slide with SlideID 258 (position 4): move the "<XML/>" big text
#    and the "sbJSON" logo group up/left on the Reference Items slide.

$p = $ppt.ActivePresentation

# Small helper: EMU -> points, with a hair of epsilon so that the
# COM layer's float rounding lands on the exact target EMU value
# instead of occasionally truncating one unit short.
function EMUToPt($emu) {
    return ($emu / 12700) + 0.00001
}

# ---- Slide with SlideID 260 (position 1) ----------------------------
$slide260 = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 260) {
        $slide260 = $p.Slides.Item($i)
    }
}

$harvestsText = $slide260.Shapes.Item("TextBox 68")
$harvestsText.Left = EMUToPt(8298398)
$harvestsText.Top  = EMUToPt(1521927)

$harvestsPic = $slide260.Shapes.Item("Picture 69")
$harvestsPic.Left = EMUToPt(8343062)
$harvestsPic.Top  = EMUToPt(1692624)

# Add the new "Products" caption textbox (ends up as shape id 4 /
# "TextBox 3" once the two scratch textboxes below are removed again,
# matching the authored id sequence).
$scratch1 = $slide260.Shapes.AddTextbox(1, 0, 0, 1, 1)
$scratch2 = $slide260.Shapes.AddTextbox(1, 0, 0, 1, 1)
$productsBox = $slide260.Shapes.AddTextbox(1, `
    (EMUToPt(8363308)), (EMUToPt(2145514)), `
    (EMUToPt(806198)), (EMUToPt(276999)))
$scratch1.Delete()
$scratch2.Delete()

$productsBox.Fill.Visible = $false
$productsBox.TextFrame.WordWrap = -1
$productsBox.TextFrame.TextRange.Text = "Products"
$productsBox.TextFrame.TextRange.Font.Name = "Arial"
$productsBox.TextFrame.TextRange.Font.NameComplexScript = "Arial"
$productsBox.TextFrame.TextRange.Font.Size = 12
$productsBox.TextFrame.AutoSize = 1

# ---- Slide with SlideID 258 (position 4) -----------------------------
$slide258 = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 258) {
        $slide258 = $p.Slides.Item($i)
    }
}

$xmlText = $slide258.Shapes.Item("TextBox 6")
$xmlText.Left = EMUToPt(3826874)
$xmlText.Top  = EMUToPt(3927435)

$sbjsonGroup = $slide258.Shapes.Item("Group 9")
$sbjsonGroup.Left = EMUToPt(3140099)
$sbjsonGroup.Top  = EMUToPt(2348892)
